$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "Login"

$ws.Range("A1").Value = "Email"
$ws.Range("B1").Value = "Password"

$ws.Range("A2").Value = "ketif98663@ikowat.com"
$ws.Range("B2").Value = "ketif98663"

$ws.Range("A3").Value = "ketif9866@ikowat.com"
$ws.Range("B3").Value = "ketif98664"

$ws.Range("A4").Value = "ketif983@ikowat.com"
$ws.Range("B4").Value = "ketif98665"

$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:ketif98663@ikowat.com")
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:ketif9866@ikowat.com")
$ws.Hyperlinks.Add($ws.Range("A4"), "mailto:ketif983@ikowat.com")

$ws.Range("A1:B4").Borders.LineStyle = 1
$ws.Range("A2:B4").Interior.Color = 65535

$ws.Columns.Item(1).AutoFit()
$ws.Columns.Item(2).AutoFit()
